$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update PM Property ID values in column B (rows 2-10)
$ws.Range("B2").Value = 22178843
$ws.Range("B3").Value = 22178844
$ws.Range("B4").Value = 22178845
$ws.Range("B5").Value = 22178846
$ws.Range("B6").Value = 22178847
$ws.Range("B7").Value = 22178848
$ws.Range("B8").Value = 22178849
$ws.Range("B9").Value = 22178850
$ws.Range("B10").Value = 22178851

# Update the selection on the sheet
$ws.Range("B2:B10").Select()
